$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: rename ZIP -> ZIP_NEW, ZIP Code of Instiution -> ZIP Code of Institution ---
$ws.Range("A6").Value = "ZIP_NEW"
$ws.Range("B6").Value = "ZIP Code of Institution"

# --- New rows 39-42: Avg_enroll_* short names ---
$ws.Range("A39").Value = "Avg_enroll_W_12-17"
$ws.Range("A40").Value = "Avg_enroll_B_12-18"
$ws.Range("A41").Value = "Avg_enroll_H_12-19"
$ws.Range("A42").Value = "Avg_enroll_A_12-20"

# --- Descriptions for white / Hispanic / Asian (black description filled in later) ---
$ws.Range("B39").Value = "Average enrollment of white students in Texas higher educational institutions between 2012 and 2017 (%)"
$ws.Range("B41").Value = "Average enrollment of Hispanic students in Texas higher educational institutions between 2012 and 2017 (%)"
$ws.Range("B42").Value = "Average enrollment of Asian students in Texas higher educational institutions between 2012 and 2017 (%)"

# --- New rows 44-48: uni/lat/lng/zips/avg_povrate short names ---
$ws.Range("A44").Value = "uni"
$ws.Range("A45").Value = "lat"
$ws.Range("A46").Value = "lng"
$ws.Range("A47").Value = "zips"
$ws.Range("A48").Value = "avg_povrate"

# --- Row 43: Avg_enroll_12-17 (all students) ---
$ws.Range("A43").Value = "Avg_enroll_12-17"
$ws.Range("B43").Value = "Average enrollment of all students in Texas higher educational institutions between 2012 and 2017 (%)"

# --- Remaining descriptions filled in afterward ---
$ws.Range("B40").Value = "Average enrollment of black students in Texas higher educational institutions between 2012 and 2017 (%)"
$ws.Range("B47").Value = "List of zipcodes which are within 6 miles of the instution's latitude and longitude coordinates. "
$ws.Range("B48").Value = "Average poverty rate of the variable ""zips"" between 2012 and 2017"

# --- Reused descriptions for 44-46 (already existing shared strings) ---
$ws.Range("B44").Value = "Institution Name "
$ws.Range("B45").Value = "Latitude of Institution"
$ws.Range("B46").Value = "Longitude of Institution"

# Apply the left/top alignment style to A44:A48
$styleRange = $ws.Range("A44:A48")
$styleRange.HorizontalAlignment = -4131
$styleRange.VerticalAlignment = -4160

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 26.17
$ws.Columns.Item(2).ColumnWidth = 88.98

# --- Page setup (orientation) ---
$ws.PageSetup.Orientation = 1

# --- Selection / view state ---
$ws.Range("B50").Select() | Out-Null
